$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@('Achraf','Othman','Mada Center','Qatar','TRCtFdMAAAAJ','M',1985,'Informatique, Mathématiques et Ingénierie')
    ,@('Afef','Khalil','Université de La Manouba','Tunisie','XKJ-HTQAAAAJ','F',1986,'Sciences économiques, Finance et Gestion')
    ,@('Amal','Ben Cheikh','Avicenne Private Business School','Tunisie','uGv8N1kAAAAJ','F',1989,'Sciences économiques, Finance et Gestion')
    ,@('Bochra','Rabbouch','Université de Tunis','Tunisie','rxq9X58AAAAJ','F',1987,'Informatique, Mathématiques et Ingénierie')
    ,@('Hana','Rabbouch','Université de Sousse','Tunisie','yW86SmIAAAAJ','F',1988,'Informatique, Mathématiques et Ingénierie')
    ,@('Ines','Amara','King Faisal University','Arabie Saoudite','ifRBvikAAAAJ','F',1986,'Sciences économiques, Finance et Gestion')
    ,@('Mehdi','Ghommem','American University of Sharjah','UAE','Cokg8GAAAAAJ','M',1983,'Informatique, Mathématiques et Ingénierie')
    ,@('Ismail','Ktata','Université de Gabès','Tunisie','0V9BV_sAAAAJ','M',1983,'Informatique, Mathématiques et Ingénierie')
    ,@('Intissar','Moussa','Université de Sousse','Tunisie','LTv022EAAAAJ','F',1988,'Informatique, Mathématiques et Ingénierie')
    ,@('Mohamed A.','Bahloul','Alfaisal University','Saudi Arabia','GTIvdXUAAAAJ','M',1991,'Informatique, Mathématiques et Ingénierie')
    ,@('Azza','Mensi','Canadian Institutes of Health Research','Canada','nKvcr2gAAAAJ','F',1985,'Médecine, Biologie et Sciences de la Santé')
    ,@('Imed','Madhi','Université Centrale','Tunisie','XVUtQ74AAAAJ','M',1983,'Physique et Astronomie')
    ,@('Fahmi','Alila','Université de Nantes','France','1v1t5G4AAAAJ','M',1990,'Informatique, Mathématiques et Ingénierie')
    ,@('Amal','Nammouchi','Karlstad University','Sweden','CALEjIEAAAAJ','F',1996,'Informatique, Mathématiques et Ingénierie')
    ,@('Sabrine','Ziri','Institut supérieur de l''aéronautique et de l''espace','France','sdhIDjYAAAAJ','F',1992,'Chimie et Sciences des Matériaux')
)

$startRow = 49
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Match the "Genre" column (F) cell style used elsewhere in the sheet (Arial 8pt)
$ws.Range("F2").Copy()
$ws.Range("F49:F63").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the view/selection state recorded after the edit
[void]$ws.Range("I74").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
